# "Catch and survival calculations... almost complete"
#
# - Rename the (empty) "GrowthChange" sheet to "Recruitment" and fill in the
#   steepness / recruitment-deviation parameters.
# - Update the Control sheet's selection to cover its whole used range.
# - Make Recruitment the tab-selected / active sheet (it was Growth before).

$wb = $excel.ActiveWorkbook

# --- Rename GrowthChange -> Recruitment and populate its parameter table ---
$ws = $wb.Worksheets.Item("GrowthChange")
$ws.Name = "Recruitment"

$ws.Range("A1").Value = "Par"
$ws.Range("B1").Value = "Value"
$ws.Range("C1").Value = "Notes"

$ws.Range("A2").Value = "h"
$ws.Range("B2").Value = 0.75
$ws.Range("C2").Value = "Steepness"

$ws.Range("A3").Value = "sigma_rec"
$ws.Range("B3").Value = 0.9
$ws.Range("C3").Value = "Lognormal recruitment standard deviation (0.55-0.9)"

# --- Control sheet: selection now spans the full used range A1:C3 ---
$ctrl = $wb.Worksheets.Item("Control")
$ctrl.Activate()
[void]$ctrl.Range("A1:C3").Select()

# --- Growth sheet: no longer the tab-selected sheet ---
$growth = $wb.Worksheets.Item("Growth")
[void]$growth.Select()

# --- Recruitment becomes the active / tab-selected sheet, cursor at E15 ---
$ws.Activate()
[void]$ws.Range("E15").Select()

$wb.Save()
